# CS320-Sp19-Teams.xlsx — "updated tutors and teams"
# Renames "Team 3" to its real project title, fills in all of the
# team rosters that were previously left as TBD, highlights a few
# teams in yellow, and restores the last-used selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------
# Section 101 table (rows 6-10): Team 1 / Team 2 / Team 3 / Team 4
# ---------------------------------------------------------------
$ws.Range("D6").Value = "Team 3: Physical Model Website and DB"

$ws.Range("B7").Value = "Jason Bady"
$ws.Range("C7").Value = "Bill Abram"
$ws.Range("D7").Value = "Tyler Kautz"
$ws.Range("E7").Value = "Lucas Gartrell"

$ws.Range("B8").Value = "Quintin Herb"
$ws.Range("C8").Value = "Dennis Chism"
$ws.Range("D8").Value = "Jake Stambaugh"
$ws.Range("E8").Value = "Cole Rohrbaugh"

$ws.Range("B9").Value = "Patrick Nelson"
$ws.Range("C9").Value = "Alyssa Gross"
$ws.Range("D9").Value = "Trevor Swan"
$ws.Range("E9").Value = "Ben Yanick"

$ws.Range("B10").Value = "Will Wyatt"
$ws.Range("C10").Value = "Joshua Grove"

# Highlight the Team 2 column (roster filled in / confirmed)
$ws.Range("C7:C10").Interior.Color = 65535

# ---------------------------------------------------------------
# Section 102 table (rows 14-18): Team 5 / Team 6 / Team 7 / Team 8
# ---------------------------------------------------------------
$ws.Range("B15").Value = "Collin Brandt"
$ws.Range("C15").Value = "Tom Herbine"
$ws.Range("D15").Value = "Sam Cesrario"
$ws.Range("E15").Value = "Alex Doyle"

$ws.Range("B16").Value = "Tim Jefferson"
$ws.Range("C16").Value = "Ethan Hostle"
$ws.Range("D16").Value = "Kyle Leatherman"
$ws.Range("E16").Value = "Vince Maresca"

$ws.Range("B17").Value = "Tom McAdams"
$ws.Range("C17").Value = "Mike Reinhart"
$ws.Range("D17").Value = "Josiah Sam"
$ws.Range("E17").Value = "Jennifer Rhine"

$ws.Range("E18").Value = "Duncab Smucker"
$ws.Range("E18").Font.Bold = $true

# The whole Team 5-8 block is fully staffed - highlight it all
$ws.Range("B15:E18").Interior.Color = 65535

# ---------------------------------------------------------------
# Section 103 table (rows 22-26): Team 9 / Team 10 / Team 11 / Team 12
# ---------------------------------------------------------------
$ws.Range("B23").Value = "Adrian Castro"
$ws.Range("C23").Value = "Jaden Marini"
$ws.Range("D23").Value = "Dakota Hilbert"
$ws.Range("E23").Value = "Jason Steinberg"

$ws.Range("B24").Value = "Darnell Hill"
$ws.Range("C24").Value = "John Steck"
$ws.Range("D24").Value = "Zack Ronayne"
$ws.Range("E24").Value = "Bill Taylor"

$ws.Range("B25").Value = "Zack Hirs"
$ws.Range("C25").Value = "Chase Traina"
$ws.Range("D25").Value = "Madison Tibbett"
$ws.Range("E25").Value = "Travis Wetzel"

$ws.Range("B26").Value = "Joe Landau"

# Highlight the Team 11 / Team 12 columns
$ws.Range("D23:E26").Interior.Color = 65535

# ---------------------------------------------------------------
# Restore the sheet's last active selection
# ---------------------------------------------------------------
$ws.Range("C11").Select()
